# Regenerate save_data column G ("K") with newly calculated values.
# (commit: regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G (rows 2-31), per the recalculated "K" series.
$newK = @{
    2  = 5
    3  = 4
    4  = 3
    5  = 9
    6  = 5
    7  = 6
    8  = 3
    9  = 6
    10 = 7
    11 = 8
    12 = 4
    13 = 6
    14 = 1
    15 = 4
    16 = 5
    17 = 5
    18 = 4
    19 = 6
    20 = 4
    21 = 5
    22 = 4
    23 = 7
    24 = 3
    25 = 5
    26 = 5
    27 = 6
    28 = 2
    29 = 4
    30 = 3
    31 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
